$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text / translation fixes (shared strings) ---
# Fix typo in the Arabic word for "Point Reference" column header (B1):
# النوصيل -> التوصيل (keep the original non-breaking-space spacing intact)
$ws.Range("B1").Value = "Point Reference" + [char]0xA0 + " " + [char]0xA0 + " " + [char]0xA0 + "رقم نقطة التوصيل"

# Add a required-field marker " *" to the "Description" column header (I1)
$ws.Range("I1").Value = "Description *     الوصف"

# --- Column width change for column C ---
$ws.Columns.Item(3).ColumnWidth = 38.43

# --- View state: active cell / selection ---
$ws.Activate()
$ws.Range("J1").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$win.ScrollRow = 1

# --- Data validation update for column C ---
# Previously a list validation ("Household electronics,Mobiles"); now just a
# plain allow-blank validation with no restriction list.
$ws.Range("C1:C1048576").Validation.Modify(0, 1, 1, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value)
